$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")

# Widen column B (was 12.1640625 -> now 19 characters)
$ws.Columns.Item(2).ColumnWidth = 18.285714285714285

# New custom width on column E, reusing the old column B width (12.1640625 characters)
$ws.Columns.Item(5).ColumnWidth = 11.428571428571429

# F16: parent netNode E1 -> E2
$ws.Range("F16").Value = "E2"

# F17: parent netNode E2 -> E3
$ws.Range("F17").Value = "E3"

# Move the active selection from G5 to I15
$ws.Range("I15").Select()
